$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("git")
Write-Host $ws.Name
Write-Host $ws.Range("B2").Text
